# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) to the classification results sheet and
# updates the refit Prediction/Error/Cross-Entropy values in columns D/E/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Label" header in H1, matching the existing header style ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Populate the new "Label" column (H2:H21) ---
$labels = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
$r = 2
foreach ($lbl in $labels) {
    $ws.Cells.Item($r, 8).Value = $lbl
    $r++
}

# --- Update refitted Prediction / Error / Cross Entropy Loss values ---
$ws.Range("D2").Value = 0.7265774850554008
$ws.Range("E2").Value = 0.7265774850554008

$ws.Range("D3").Value = 0.3873500000952162
$ws.Range("E3").Value = 0.3873500000952162

$ws.Range("D4").Value = 0.577848931864977
$ws.Range("E4").Value = 0.577848931864977

$ws.Range("D5").Value = 0.5196081268459193
$ws.Range("E5").Value = 0.5196081268459193

$ws.Range("D6").Value = 0.380469579434884
$ws.Range("E6").Value = 0.380469579434884

$ws.Range("D7").Value = 0.3876217229866097
$ws.Range("E7").Value = 0.6123782770133903

$ws.Range("D8").Value = 0.6206681097185163
$ws.Range("E8").Value = 0.3793318902814837

$ws.Range("D10").Value = 0.7650655720243033
$ws.Range("E10").Value = 0.2349344279756967

$ws.Range("F11").Value = 0.6207945346832275
